$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking strings
# are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.526.40"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").Value = "1.918.45"
$ws.Range("E3").Value = "  -0.17%  "

$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "245.72"
$ws.Range("E5").Value = "  +1.11%  "

$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "0.4795"
$ws.Range("E7").Value = "  +1.75%  "

$ws.Range("D8").Value = "0.2899"
$ws.Range("E8").Value = "  +0.80%  "

$ws.Range("D9").Value = "0.06722"
$ws.Range("E9").Value = "  -0.46%  "

$ws.Range("D10").Value = "110.99"
$ws.Range("E10").Value = "  +4.18%  "

$ws.Range("D11").Value = "19.05"
$ws.Range("E11").Value = "  +4.16%  "

$ws.Range("D12").Value = "1.910.78"
$ws.Range("E12").Value = "  -0.53%  "

$ws.Range("D13").Value = "0.07568"
$ws.Range("E13").Value = "  -2.52%  "

$ws.Range("D14").Value = "5.274"
$ws.Range("E14").Value = "  -0.49%  "

$ws.Range("D15").Value = "0.6679"
$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("D16").Value = "299.75"
$ws.Range("E16").Value = "  +2.59%  "

$ws.Range("D17").Value = "30.523.45"
$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "5.641"
$ws.Range("E18").Value = "  +6.61%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "12.99"
$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").Value = "0.000007578"
$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").Value = "2.165.36"
$ws.Range("E22").Value = "  +0.58%  "

$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").Value = "6.491"
$ws.Range("E24").Value = "  +4.53%  "

$ws.Range("D25").Value = "9.483"
$ws.Range("E25").Value = "  +1.23%  "

$ws.Range("D26").Value = "164.66"
$ws.Range("E26").Value = "  -2.16%  "

$ws.Range("D27").Value = "20.29"
$ws.Range("E27").Value = "  -4.78%  "

$ws.Range("D28").Value = "2.111"
$ws.Range("E28").Value = "  +0.63%  "

$ws.Range("D29").Value = "0.1077"
$ws.Range("E29").Value = "  +0.67%  "

$ws.Range("E30").Value = "  +2.34%  "

$ws.Range("D31").Value = "4.170"
$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("D32").Value = "4.049"
$ws.Range("E32").Value = "  +1.16%  "

$ws.Range("D33").Value = "0.05002"
$ws.Range("E33").Value = "  -0.68%  "

$ws.Range("D34").Value = "0.7372"
$ws.Range("E34").Value = "  -0.86%  "

$ws.Range("D35").Value = "1.137"
$ws.Range("E35").Value = "  -1.47%  "

$ws.Range("D37").Value = "2.722"
$ws.Range("E37").Value = "  -0.23%  "

$ws.Range("D38").Value = "0.02036"
$ws.Range("E38").Value = "  -3.72%  "

$ws.Range("E39").Value = "  +0.07%  "

$ws.Range("D40").Value = "111.08"
$ws.Range("E40").Value = "  +0.75%  "

$ws.Range("D41").Value = "2.023"
$ws.Range("E41").Value = "  -2.47%  "

$ws.Range("D42").Value = "0.4437"
$ws.Range("E42").Value = "  +3.92%  "

$ws.Range("D43").Value = "72.43"
$ws.Range("E43").Value = "  +7.28%  "

$ws.Range("D44").Value = "0.8625"

$ws.Range("D45").Value = "5.876"
$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("D46").Value = "1.0000"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").Value = "49.38"
$ws.Range("E47").Value = "  -0.64%  "

$ws.Range("D48").Value = "7.283"
$ws.Range("E48").Value = "  +1.26%  "

$ws.Range("D49").Value = "9.311"
$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").Value = "0.2549"
$ws.Range("E50").Value = "  +4.08%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.1230"
$ws.Range("E51").Value = "  +0.91%  "

# Restore default (no explicit) style on column D now that values are set,
# so the cell style matches the original workbook (no style index).
$ws.Range("D2:D51").Style = "Normal"
